# Scottish Module Input files modification
# - Rename "GeneralTaxRateMonthly" sheet/content to "GeneralTaxRateWeekly"
# - Rename "ProcessPayrollForMonthlyTax" sheet/content to "ProcessPayrollForWeeklyTax"
# - Update "DO NOT TOUCH AUTOMATION EMP 105" -> "DO NOT TOUCH AUTOMATION EMP 107"
#   on the GeneralTaxRate*, ProcessPayrollFor*Tax and TestReports sheets.

$wb = $excel.ActiveWorkbook

$wsFirst   = $wb.Worksheets.Item(1)   # "first"
$wsMonthly = $wb.Worksheets.Item(2)   # "GeneralTaxRateMonthly"
$wsProcess = $wb.Worksheets.Item(3)   # "ProcessPayrollForMonthlyTax"
$wsReports = $wb.Worksheets.Item(4)   # "TestReports"

# --- Rename sheet tabs -------------------------------------------------
$wsMonthly.Name = "GeneralTaxRateWeekly"
$wsProcess.Name = "ProcessPayrollForWeeklyTax"

# --- Update cell text that referenced the old sheet names on "first" ---
$wsFirst.Range("A3").Value = "GeneralTaxRateWeekly"
$wsFirst.Range("A4").Value = "ProcessPayrollForWeeklyTax"

# --- Update the "DO NOT TOUCH AUTOMATION EMP 105" -> 107 markers -------
$wsMonthly.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

# --- Update selection / active sheet to match the saved view state -----
[void]$wsMonthly.Range("A2").Select()
[void]$wsProcess.Range("B2").Select()
[void]$wsReports.Range("B2").Select()

[void]$wsFirst.Activate()
[void]$wsFirst.Range("A3").Select()
